$wb = $excel.ActiveWorkbook

# Sheets "展览" and "全部类型" both contain the same table with updated counts
# in column F ("想去人数") for rows 2, 4 and 5.
$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    $ws.Range("F2").Value = 9703
    $ws.Range("F4").Value = 34
    $ws.Range("F5").Value = 550
}
